# Append a new data row to the "OrderEntry" sheet, mirroring row 2 but with
# a new Product value ("UPC_ANYDAY_1(B)"), same TM_ID ("S1002") and same
# DelToDate (21-Mar-2023), then move the selection to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrderEntry")

$ws.Range("A3").Value = "UPC_ANYDAY_1(B)"
$ws.Range("B3").Value = "S1002"
$ws.Range("C3").Value2 = $ws.Range("C2").Value2
$ws.Range("C3").NumberFormat = $ws.Range("C2").NumberFormat

$ws.Activate() | Out-Null
$ws.Range("B3").Select() | Out-Null
